$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-17 03:34:01"
$wsZh.Range("D3").Value = "2016-02-17 03:34:01"
$wsZh.Range("G2").Value = "2016-02-17 03:34:46"
$wsZh.Range("G3").Value = "2016-02-17 03:34:46"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-17 03:34:11"
$wsDe.Range("D3").Value = "2016-02-17 03:34:11"
$wsDe.Range("G2").Value = "2016-02-17 03:35:03"
$wsDe.Range("G3").Value = "2016-02-17 03:35:03"
